# Updates cryptos list values per upstream data refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (prices with "." thousands separators, percentages, names, links) ---
$ws.Range("D2").Value = "42.962.72"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.387.31"
$ws.Range("E3").Value = "  +6.18%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("E5").Value = "  +10.98%  "
$ws.Range("E6").Value = "  -7.14%  "
$ws.Range("E7").Value = "  +3.07%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +7.43%  "
$ws.Range("E10").Value = "  -5.24%  "
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("E12").Value = "  -3.59%  "
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("E14").Value = "  +13.07%  "
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "2.749.94"
$ws.Range("E16").Value = "  +6.43%  "
$ws.Range("D17").Value = "2.380.84"
$ws.Range("E17").Value = "  +3.92%  "
$ws.Range("D18").Value = "43.097.44"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("E19").Value = "  +8.76%  "
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("E22").Value = "  +5.37%  "
$ws.Range("E23").Value = "  +8.29%  "
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("E25").Value = "  +7.71%  "
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  +4.07%  "
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +4.52%  "
$ws.Range("E34").Value = "  +3.75%  "
$ws.Range("E35").Value = "  +5.04%  "
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  -2.55%  "
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("E40").Value = "  +15.90%  "
$ws.Range("E41").Value = "  +19.57%  "
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("E44").Value = "  +14.33%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("E47").Value = "  +49.24%  "
$ws.Range("E48").Value = "  +8.24%  "
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("E51").Value = "  +11.57%  "

# --- Numeric-looking price cells: force Text format first so Excel keeps them as strings ---
# (matches the source data which stores every Price cell as text, not a number)
$numericCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D30", "D31", "D32", "D34", "D36", "D37", "D39", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "329.40"
$ws.Range("D6").Value = "104.22"
$ws.Range("D7").Value = "0.648"
$ws.Range("D9").Value = "0.652"
$ws.Range("D10").Value = "41.72"
$ws.Range("D11").Value = "0.0939"
$ws.Range("D12").Value = "8.68"
$ws.Range("D14").Value = "17.23"
$ws.Range("D19").Value = "7.83"
$ws.Range("D20").Value = "0.0000109"
$ws.Range("D21").Value = "76.63"
$ws.Range("D22").Value = "3.68"
$ws.Range("D23").Value = "272.55"
$ws.Range("D24").Value = "2.42"
$ws.Range("D25").Value = "9.67"
$ws.Range("D26").Value = "11.74"
$ws.Range("D28").Value = "23.00"
$ws.Range("D30").Value = "175.04"
$ws.Range("D31").Value = "37.55"
$ws.Range("D32").Value = "3.17"
$ws.Range("D34").Value = "5.89"
$ws.Range("D36").Value = "4.96"
$ws.Range("D37").Value = "4.17"
$ws.Range("D39").Value = "0.107"
$ws.Range("D41").Value = "1.59"
$ws.Range("D42").Value = "0.234"
$ws.Range("D43").Value = "69.78"
$ws.Range("D44").Value = "121.72"
$ws.Range("D46").Value = "12.40"
$ws.Range("D47").Value = "90.45"
$ws.Range("D48").Value = "9.36"
$ws.Range("D49").Value = "5.53"
$ws.Range("D50").Value = "1.32"
$ws.Range("D51").Value = "0.489"

# Reset style index on those cells back to the default "Normal" so no stray text-format style
# is left behind on cells that originally had none.
foreach ($addr in $numericCells) {
    $ws.Range($addr).Style = "Normal"
}
